$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '49.395.23'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.45%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.628.41'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.39%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.18%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '112.46'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.82%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '323.53'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.527'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -1.36%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.13%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.543'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.38%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '39.86'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.78'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -4.36%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0812'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.19%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.126'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.18%  '

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.06%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.041.78'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.642.23'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.86%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.859'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.01%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '49.360.46'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.02'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -1.70%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.91'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -3.26%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.70'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -1.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0945'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.45%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '269.95'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -3.38%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '68.71'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -5.72%  '

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.21'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.62%  '

$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.07%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.31'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +3.65%  '

$ws.Range("B29").Value = 'Toncoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.22'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.48%  '

$ws.Range("B30").Value = 'InjectiveProtocol'
$ws.Range("C30").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '35.02'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.07%  '

$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.138'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.71%  '

$ws.Range("B32").Value = 'OKB'
$ws.Range("C32").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '49.58'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.60%  '

$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.49'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.82%  '

$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0813'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.51%  '

$ws.Range("B35").Value = 'FirstDigitalUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("B36").Value = 'Celestia'
$ws.Range("C36").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.94'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -4.55%  '

$ws.Range("B37").Value = 'RenderToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.91'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.98%  '

$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.04'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.18%  '

$ws.Range("B39").Value = 'LidoDAOToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.13'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.55%  '

$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '127.26'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +3.20%  '

$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.111'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.64%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '22.20'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.23%  '

$ws.Range("B43").Value = 'VeChain'
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0322'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.49%  '

$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.14'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -3.91%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.059.49'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.05%  '

$ws.Range("B46").Value = 'Stacks'
$ws.Range("C46").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.15'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +6.79%  '

$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.23'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.26%  '

$ws.Range("B48").Value = 'ApeXProtocol'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.13'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -8.31%  '

$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.93'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.43%  '

$ws.Range("B50").Value = 'MultiversX'
$ws.Range("C50").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '59.11'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +1.18%  '

$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.21'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.15%  '
